# Auto-generated Excel COM-interop script to apply scheduled-runner market data refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the Ultima Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1614.2858
$ws.Range("J46").Value = 1716.6666
$ws.Range("L46").Value = 5149.9998
$ws.Range("N46").Value = -5387.9998

$ws.Range("H60").Value = 1614.2858
$ws.Range("J60").Value = 1716.6666
$ws.Range("L60").Value = 5149.9998
$ws.Range("N60").Value = -6117.9998

$ws.Range("H132").Value = 2407.6309
$ws.Range("I132").Value = 1973.0518
$ws.Range("J132").Value = 6008.4287
$ws.Range("K132").Value = 5919.1554
$ws.Range("L132").Value = 18025.2861
$ws.Range("M132").Value = -3389.1554
$ws.Range("N132").Value = -23085.2861

$ws.Range("H137").Value = 5406267.5
$ws.Range("I137").Value = 832.381
$ws.Range("J137").Value = 12500901
$ws.Range("K137").Value = 2497.143
$ws.Range("L137").Value = 37502703
$ws.Range("M137").Value = 52.85699999999997
$ws.Range("N137").Value = -37507803

$ws.Range("H138").Value = 2268.4924
$ws.Range("I138").Value = 1300.5161
$ws.Range("J138").Value = 3102.0278
$ws.Range("K138").Value = 3901.5483
$ws.Range("L138").Value = 9306.0834
$ws.Range("M138").Value = 1238.4517
$ws.Range("N138").Value = -19586.0834

$ws.Range("H141").Value = 1202.55
$ws.Range("I141").Value = 888.80554
$ws.Range("J141").Value = 4026.25
$ws.Range("K141").Value = 2666.41662
$ws.Range("L141").Value = 12078.75
$ws.Range("M141").Value = 2513.58338
$ws.Range("N141").Value = -22438.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1351.0714
$ws.Range("I2").Value = 1031.7142
$ws.Range("J2").Value = 1670.4286
$ws.Range("K2").Value = 1031.7142
$ws.Range("L2").Value = 1670.4286
$ws.Range("M2").Value = -918.7141999999999
$ws.Range("N2").Value = -1896.4286

$ws.Range("H32").Value = 7110.1626
$ws.Range("I32").Value = 6727.575
$ws.Range("J32").Value = 11100
$ws.Range("K32").Value = 6727.575
$ws.Range("L32").Value = 11100
$ws.Range("M32").Value = -6440.575
$ws.Range("N32").Value = -11674

$ws.Range("H45").Value = 3862.4
$ws.Range("I45").Value = 6006
$ws.Range("J45").Value = 2433.3333
$ws.Range("K45").Value = 6006
$ws.Range("L45").Value = 2433.3333
$ws.Range("M45").Value = -5629
$ws.Range("N45").Value = -3187.3333

$ws.Range("H61").Value = 7247483
$ws.Range("I61").Value = 9435002
$ws.Range("J61").Value = 1325.9375
$ws.Range("K61").Value = 9435002
$ws.Range("L61").Value = 1325.9375
$ws.Range("M61").Value = -9434790
$ws.Range("N61").Value = -1749.9375

$ws.Range("H110").Value = 1376.375
$ws.Range("I110").Value = 1376.375
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1376.375
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 668.625

$ws.Range("H116").Value = 1351.0714
$ws.Range("I116").Value = 1031.7142
$ws.Range("J116").Value = 1670.4286
$ws.Range("K116").Value = 1031.7142
$ws.Range("L116").Value = 1670.4286
$ws.Range("M116").Value = 1262.2858
$ws.Range("N116").Value = -6258.4286

$ws.Range("H122").Value = 2836.8333
$ws.Range("I122").Value = 2498.6667
$ws.Range("J122").Value = 3175
$ws.Range("K122").Value = 7496.000100000001
$ws.Range("L122").Value = 9525
$ws.Range("M122").Value = -5046.000100000001
$ws.Range("N122").Value = -14425

$ws.Range("H132").Value = 8066477
$ws.Range("I132").Value = 10871313
$ws.Range("J132").Value = 2573
$ws.Range("K132").Value = 32613939
$ws.Range("L132").Value = 7719
$ws.Range("M132").Value = -32611409
$ws.Range("N132").Value = -12779

$ws.Range("H136").Value = 7247483
$ws.Range("I136").Value = 9435002
$ws.Range("J136").Value = 1325.9375
$ws.Range("K136").Value = 28305006
$ws.Range("L136").Value = 3977.8125
$ws.Range("M136").Value = -28302456
$ws.Range("N136").Value = -9077.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1351.0714
$ws.Range("I3").Value = 1031.7142
$ws.Range("J3").Value = 1670.4286
$ws.Range("K3").Value = 1031.7142
$ws.Range("L3").Value = 1670.4286
$ws.Range("M3").Value = -917.7141999999999
$ws.Range("N3").Value = -1898.4286

$ws.Range("H134").Value = 2075.5334
$ws.Range("I134").Value = 1310.2565
$ws.Range("J134").Value = 7049.8335
$ws.Range("K134").Value = 3930.7695
$ws.Range("L134").Value = 21149.5005
$ws.Range("M134").Value = -1395.7695
$ws.Range("N134").Value = -26219.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5379605
$ws.Range("I31").Value = 3736.7334
$ws.Range("J31").Value = 19609844
$ws.Range("K31").Value = 3736.7334
$ws.Range("L31").Value = 19609844
$ws.Range("M31").Value = -3441.7334
$ws.Range("N31").Value = -19610434

$ws.Range("H34").Value = 5379605
$ws.Range("I34").Value = 3736.7334
$ws.Range("J34").Value = 19609844
$ws.Range("K34").Value = 3736.7334
$ws.Range("L34").Value = 19609844
$ws.Range("M34").Value = -3534.7334
$ws.Range("N34").Value = -19610248

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0

$ws.Range("H88").Value = 50343
$ws.Range("J88").Value = 50343
$ws.Range("L88").Value = 50343
$ws.Range("N88").Value = -51155

$ws.Range("H91").Value = 50343
$ws.Range("J91").Value = 50343
$ws.Range("L91").Value = 50343
$ws.Range("N91").Value = -53151

$ws.Range("H134").Value = 1429.5278
$ws.Range("I134").Value = 1527.16
$ws.Range("J134").Value = 1207.6364
$ws.Range("K134").Value = 4581.48
$ws.Range("L134").Value = 3622.9092
$ws.Range("M134").Value = -2046.48
$ws.Range("N134").Value = -8692.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 538.75
$ws.Range("I60").Value = 430
$ws.Range("K60").Value = 1290
$ws.Range("M60").Value = -1039

$ws.Range("H63").Value = 6378.1665
$ws.Range("I63").Value = 5233
$ws.Range("K63").Value = 15699
$ws.Range("M63").Value = -14950

$ws.Range("H66").Value = 6378.1665
$ws.Range("I66").Value = 5233
$ws.Range("K66").Value = 47097
$ws.Range("M66").Value = -43353

$ws.Range("H131").Value = 847.64
$ws.Range("I131").Value = 562.5
$ws.Range("J131").Value = 859.5208
$ws.Range("K131").Value = 1687.5
$ws.Range("L131").Value = 2578.5624
$ws.Range("M131").Value = 3352.5
$ws.Range("N131").Value = -12658.5624

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 9000
$ws.Range("I53").Value = 10000
$ws.Range("J53").Value = 8000
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 8000
$ws.Range("M53").Value = -9369
$ws.Range("N53").Value = -9262

$ws.Range("H102").Value = 1903.1034
$ws.Range("I102").Value = 1727.6
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1727.6
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -105.5999999999999
$ws.Range("N102").Value = -6244

$ws.Range("H122").Value = 2669714.2
$ws.Range("I122").Value = 3511203.5
$ws.Range("J122").Value = 4997.8335
$ws.Range("K122").Value = 10533610.5
$ws.Range("L122").Value = 14993.5005
$ws.Range("M122").Value = -10531160.5
$ws.Range("N122").Value = -19893.5005

$ws.Range("H126").Value = 3832.5527
$ws.Range("I126").Value = 2550.2632
$ws.Range("J126").Value = 5114.8423
$ws.Range("K126").Value = 7650.7896
$ws.Range("L126").Value = 15344.5269
$ws.Range("M126").Value = -5180.7896
$ws.Range("N126").Value = -20284.5269

$ws.Range("H132").Value = 2963.0327
$ws.Range("I132").Value = 2090.3262
$ws.Range("J132").Value = 5639.3335
$ws.Range("K132").Value = 6270.9786
$ws.Range("L132").Value = 16918.0005
$ws.Range("M132").Value = -3740.9786
$ws.Range("N132").Value = -21978.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6369
$ws.Range("I40").Value = 7850
$ws.Range("J40").Value = 4888
$ws.Range("K40").Value = 7850
$ws.Range("L40").Value = 4888
$ws.Range("M40").Value = -7714
$ws.Range("N40").Value = -5160

$ws.Range("H122").Value = 6252.6665
$ws.Range("I122").Value = 6822.5454
$ws.Range("J122").Value = 5357.143
$ws.Range("K122").Value = 20467.6362
$ws.Range("L122").Value = 16071.429
$ws.Range("M122").Value = -18017.6362
$ws.Range("N122").Value = -20971.429

$ws.Range("H134").Value = 56807.96
$ws.Range("J134").Value = 56807.96
$ws.Range("L134").Value = 56807.96
$ws.Range("N134").Value = -66947.95999999999

$ws.Range("H136").Value = 9808095
$ws.Range("I136").Value = 15627021
$ws.Range("J136").Value = 7797.3687
$ws.Range("K136").Value = 46881063
$ws.Range("L136").Value = 23392.1061
$ws.Range("M136").Value = -46878513
$ws.Range("N136").Value = -28492.1061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0

$ws.Range("H122").Value = 3242.3333
$ws.Range("I122").Value = 3626
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 10878
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -8428
$ws.Range("N122").Value = -12325

$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws.Range("H132").Value = 3900.2354
$ws.Range("I132").Value = 2763.7273
$ws.Range("J132").Value = 5983.8335
$ws.Range("K132").Value = 8291.1819
$ws.Range("L132").Value = 17951.5005
$ws.Range("M132").Value = -5761.1819
$ws.Range("N132").Value = -23011.5005
